$d = $word.ActiveDocument

# 1. Remove the hidden "_GoBack" bookmark from the first (empty) paragraph.
#    Word re-creates this automatically while editing; stripping it here
#    matches the saved/"clean" document state (<w:p/> with no bookmark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

# 2. Header cell: "Full Name" -> "Name"
$null = $d.Content.Find.Execute("Full Name", $true, $false, $false, $false, $false, $true, 1, $false, "Name", 2)

# 3. Header cell: "Background (e.g. Engineering, Science, Math, CS, IT, Biz, etc)"
#    -> "Undergraduate Background (degree, major, etc)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 4)
$s = $cell.Range.Start

$oldRunA = "Background (e.g. Engineering, Science,"
$oldRunB = " Math, CS, IT, Biz, "
$newRunA = "Undergraduate "
$newRunB = "Background"
$newRunC = " (degree, major, "

$lenA = $oldRunA.Length
$lenB = $oldRunB.Length

# Insert the new " (degree, major, " text right before "etc" (i.e. right
# after the original " Math, CS, IT, Biz, " run) first, while the
# preceding offsets are still the original ones.
$insertPoint = $d.Range($s + $lenA + $lenB, $s + $lenA + $lenB)
$insertPoint.InsertBefore($newRunC)

# Replace "Background (e.g. Engineering, Science," with "Undergraduate "
$rA = $d.Range($s, $s + $lenA)
$rA.Text = $newRunA

# The former " Math, CS, IT, Biz, " run now starts (lenA - newRunA.Length)
# chars earlier, and is still $lenB chars long.
$shift = $lenA - $newRunA.Length
$rB = $d.Range($s + $lenA - $shift, $s + $lenA - $shift + $lenB)
$rB.Text = $newRunB

# 4. Resize table columns 4 & 5 (Background / Contact) - dxa 3240->2790 and 2070->2520
#    i.e. points 162->139.5 and 103.5->126
$t.Columns.Item(4).Width = 139.5
$t.Columns.Item(5).Width = 126
